# Add a new "Addressing mode(Optional)" column (O) to the node-info template.
# This reflects switching the "addressing mode" UI control from a toggle
# switcher to a select dropdown, whose two choices ("Static" / "Dynamic")
# are now documented as example values in the template, one per data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for the new column O.
$ws.Range("O1").Value = "Addressing mode(Optional)"

# Sample values for the two populated data rows (row 4 is left blank, as
# in the source change). Row 3 is written before row 2 so the shared
# strings are appended to the table in "Static" (idx 33), "Dynamic"
# (idx 34) order.
$ws.Range("O3").Value = "Static"
$ws.Range("O2").Value = "Dynamic"

# Match the new column's width to the source template.
$ws.Columns("O").ColumnWidth = 26.36

# Restore the active selection to match the saved view.
$null = $ws.Range("G17").Select()
